$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-14 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-15 Saturday", 2) | Out-Null
$d.Content.Find.Execute("16÷6=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=9, 1", 2) | Out-Null
$d.Content.Find.Execute("42÷3=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷7=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=35, 0", 2) | Out-Null
$d.Content.Find.Execute("45÷2=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=5, 1", 2) | Out-Null
$d.Content.Find.Execute("82÷9=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=21, 0", 2) | Out-Null
$d.Content.Find.Execute("57÷7=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "57÷9=6, 3", 2) | Out-Null
$d.Content.Find.Execute("95÷6=15, 5", $true, $false, $false, $false, $false, $true, 1, $false, "58÷8=7, 2", 2) | Out-Null
$d.Content.Find.Execute("73÷4=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=38, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷2=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "90÷3=30, 0", 2) | Out-Null
$d.Content.Find.Execute("47÷7=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "58÷8=7, 2", 2) | Out-Null
$d.Content.Find.Execute("11÷2=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷3=10, 0", 2) | Out-Null
$d.Content.Find.Execute("96÷2=48, 0", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=6, 5", 2) | Out-Null
$d.Content.Find.Execute("92÷4=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=9, 3", 2) | Out-Null
$d.Content.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷7=2, 4", 2) | Out-Null
$d.Content.Find.Execute("84÷8=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "57÷2=28, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷6=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 2) | Out-Null
$d.Content.Find.Execute("28÷5=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "19÷5=3, 4", 2) | Out-Null
$d.Content.Find.Execute("65÷9=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "51÷8=6, 3", 2) | Out-Null
$d.Content.Find.Execute("22÷5=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "28÷8=3, 4", 2) | Out-Null
$d.Content.Find.Execute("82÷4=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "11÷7=1, 4", 2) | Out-Null
$d.Content.Find.Execute("66÷8=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "60÷3=20, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷2=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=4, 0", 2) | Out-Null
$d.Content.Find.Execute("51÷5=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=18, 3", 2) | Out-Null
$d.Content.Find.Execute("43÷4=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "23÷4=5, 3", 2) | Out-Null
$d.Content.Find.Execute("62÷3=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "36÷5=7, 1", 2) | Out-Null
